$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (row 1) - force text format so Excel doesn't
# auto-convert these "Month Year" strings into date serial numbers
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "May 2024"

$ws.Range("G1").NumberFormat = "@"
$ws.Range("G1").Value = "June 2024"

# Update data values (row 2)
$ws.Range("A2").Value = 1.934
$ws.Range("B2").Value = -0.203
$ws.Range("C2").Value = -0.032
$ws.Range("D2").Value = 0.012
$ws.Range("E2").Value = 0.022
$ws.Range("F2").Value = -0.233
$ws.Range("G2").Value = 1.502
